$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Contenu du stage" breakdown (rows 16-23): update counts (E) and percentages (G)
# Percentages are stored as plain text (e.g. "3.03 %"), so force text format
# before assigning, otherwise Excel auto-converts the string into a numeric
# percentage value.

$ws.Range("E16").Value = 1
$ws.Range("G16").NumberFormat = "@"
$ws.Range("G16").Value = "3.03 %"

$ws.Range("E17").Value = 29
$ws.Range("G17").NumberFormat = "@"
$ws.Range("G17").Value = "87.88 %"

$ws.Range("E18").Value = 0
$ws.Range("G18").NumberFormat = "@"
$ws.Range("G18").Value = "0 %"

$ws.Range("E19").Value = 2
$ws.Range("G19").NumberFormat = "@"
$ws.Range("G19").Value = "6.06 %"

$ws.Range("E20").Value = 1
$ws.Range("G20").NumberFormat = "@"
$ws.Range("G20").Value = "3.03 %"

$ws.Range("E21").Value = 0
$ws.Range("G21").NumberFormat = "@"
$ws.Range("G21").Value = "0 %"

$ws.Range("E22").Value = 0
$ws.Range("G22").NumberFormat = "@"
$ws.Range("G22").Value = "0 %"

$ws.Range("E23").Value = 0
$ws.Range("G23").NumberFormat = "@"
$ws.Range("G23").Value = "0 %"
